$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# --- Header row (row 1) ---
$ws.Range("C1").Value = "linkText"
$ws.Range("D1").Value = "expectedUrl"

# --- Row 2 (Noticias) ---
# Leading apostrophe preserves a literal leading/trailing space text value
# and keeps/creates the "quote prefix" cell format (as in the original file).
$ws.Range("C2").Value = "' Noticias   "
$ws.Range("D2").Value = "https://es.noticias.yahoo.com/"

# --- Row 3 (Deportes) ---
$ws.Range("A3").Value = "Y"
$ws.Range("B3").Value = "YNS002"
$ws.Range("C3").Value = " Deportes   "
$ws.Range("D3").Value = "https://es.sports.yahoo.com/"

# --- Row 4 (Finanzas) ---
$ws.Range("A4").Value = "Y"
$ws.Range("B4").Value = "YNS003"
$ws.Range("C4").Value = "Finanzas"
$ws.Range("D4").Value = "https://es.finance.yahoo.com/"

# --- Row 5 (TV) ---
$ws.Range("A5").Value = "Y"
$ws.Range("B5").Value = "YNS004"
$ws.Range("C5").Value = "TV"
$ws.Range("D5").Value = "https://es.vida-estilo.yahoo.com/television/"

# Give the new rows the same centered look as the rest of the table
$ws.Range("A3:A5").Value = $ws.Range("A3:A5").Value
$ws.Range("B3:D5").HorizontalAlignment = -4108
$ws.Range("B3:D5").VerticalAlignment = -4108

# Re-assert vertical centering for the header + second row, preserving
# their existing (bold / quote-prefixed) formatting.
$ws.Range("C1:D1").VerticalAlignment = -4108
$ws.Range("C2").VerticalAlignment = -4108

# Turn the URL in D2 into a real hyperlink (creates the Hyperlink style too)
$ws.Hyperlinks.Add($ws.Range("D2"), "https://es.noticias.yahoo.com/", [Type]::Missing, [Type]::Missing, "https://es.noticias.yahoo.com/") | Out-Null
$ws.Range("D2").VerticalAlignment = -4108
$ws.Range("D2").HorizontalAlignment = -4108

# Column D is now wider to fit the urls
$ws.Columns("D").ColumnWidth = 44.16

$ws.Range("C5").Select() | Out-Null

Write-Host "edit applied"
